$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Experimental Factor Ontology (efo) version: v3.59.0 -> v3.61.0
$ws.Range("E4").Value = "v3.61.0"

# Update Disease Ontology (do) version: v2023-10-21 -> v2023-12-20
$ws.Range("E3").Value = "v2023-12-20"
